$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet right after "总计" (i.e. before "2022-Q2")
# ---------------------------------------------------------------------------
$sheetQ2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($sheetQ2)
$q3.Name = "2022-Q3"

# Header row
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q3.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$q3Data = @(
    @(0, "167301", "方正富邦中证保险主题指数（LOF）", "52.55", "93.05", "4.60", "2.4173", 5),
    @(1, "002768", "华安安进灵活配置混合A", "22.60", "74.81", "3.99", "0.9017", 6),
    @(2, "002670", "万家沪深300指数增强A", "20.85", "94.06", "2.00", "0.4170", 9),
    @(3, "519979", "长信内需成长混合A", "6.93", "92.49", "3.09", "0.2141", 9),
    @(4, "007143", "国投瑞银沪深300指数量化增强A", "9.81", "92.97", "2.12", "0.2080", 10),
    @(5, "002671", "万家沪深300指数增强C", "10.38", "94.06", "2.00", "0.2076", 9),
    @(6, "519997", "长信银利精选混合A", "4.01", "79.98", "3.58", "0.1436", 10),
    @(7, "007144", "国投瑞银沪深300指数量化增强C", "4.02", "92.97", "2.12", "0.0852", 10),
    @(8, "000877", "华泰柏瑞量化优选灵活配置混合", "6.63", "92.31", "1.10", "0.0729", 8),
    @(9, "006121", "华安双核驱动混合A", "0.83", "93.02", "4.90", "0.0407", 4),
    @(10, "257040", "国联安红利混合", "0.65", "76.59", "5.38", "0.0350", 6),
    @(11, "009059", "南方沪深300指数增强A", "1.22", "93.96", "2.60", "0.0317", 3),
    @(12, "015768", "长信内需成长混合C", "0.82", "92.49", "3.09", "0.0253", 9),
    @(13, "008184", "新华沪深300指数增强C", "1.13", "93.62", "1.68", "0.0190", 10),
    @(14, "010868", "华宝安盈混合", "5.72", "20.57", "0.33", "0.0189", 7),
    @(15, "005248", "新华沪深300指数增强A", "0.85", "93.62", "1.68", "0.0143", 10),
    @(16, "860029", "光大阳光对冲策略6个月持有期灵活配置混合C", "2.95", "63.41", "0.47", "0.0139", 7),
    @(17, "009060", "南方沪深300指数增强C", "0.39", "93.96", "2.60", "0.0101", 3),
    @(18, "002334", "汇丰晋信大盘波动精选股票A", "0.16", "85.75", "1.92", "0.0031", 9),
    @(19, "013504", "华安双核驱动混合C", "0.04", "93.02", "4.90", "0.0020", 4),
    @(20, "860028", "光大阳光对冲策略6个月持有期灵活配置混合B", "0.39", "63.41", "0.47", "0.0018", 7),
    @(21, "014572", "长信银利精选混合C", "0.05", "79.98", "3.58", "0.0018", 10),
    @(22, "006397", "长信内需成长混合E", "0.05", "92.49", "3.09", "0.0015", 9),
    @(23, "006347", "安信量化优选股票C", "0.15", "90.50", "0.75", "0.0011", 10),
    @(24, "002335", "汇丰晋信大盘波动精选股票C", "0.01", "85.75", "1.92", "0.0002", 9),
    @(25, "006346", "安信量化优选股票A", "0.03", "90.50", "0.75", "0.0002", 10),
    @(26, "860010", "光大阳光对冲策略6个月持有期灵活配置混合A", "0.03", "63.41", "0.47", "0.0001", 7),
    @(27, "016182", "华安安进灵活配置混合C", "0.00", "74.81", "3.99", "0", 6),
)

# Force columns B,C,D,E,F,G to be plain text (matches source data which stores
# these numeric-looking values as text), then fix up the one true numeric
# exception (G29) afterwards.
$lastRow = 1 + $q3Data.Length
$q3.Range("B2:G" + $lastRow).NumberFormat = "@"

for ($i = 0; $i -lt $q3Data.Length; $i++) {
    $row = $q3Data[$i]
    $r = $i + 2
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
}

# Row 29 (last data row): G column holds a real 0 number, not text "0"
$q3.Range("G" + $lastRow).NumberFormat = "General"
$q3.Range("G" + $lastRow).Value = 0

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row 2 for 2022-Q3 and
#    push the existing quarters down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# The freshly inserted row inherits the header's bold/boxed formatting -
# strip that, then copy column-A's normal numbering style onto A2.
$summary.Range("A2:D2").ClearFormats()
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 28
$summary.Cells.Item(2, 4).Value = 4.89

# Column A is a plain 0-based row counter, independent of the data that
# shifted underneath it - renumber it sequentially for all data rows.
for ($i = 0; $i -le 7; $i++) {
    $summary.Cells.Item($i + 2, 1).Value = $i
}
